$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$nl = [char]13
$tr.Text = "KÖSZÖNÖM A FIGYELMET!" + $nl + "Farkas László"

$para2 = $tr.Characters(23, 13)
$para2.Font.Italic = $true
$para2.Font.Color.RGB = 16777215
$para2.ParagraphFormat.SpaceBefore = 10
